$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B94: make it a real number instead of text "3"
$ws.Range("B94").Value = 3

# Add new row 95 with the additional annotation data
$ws.Range("A95").Value = "Ruilin"
$ws.Range("B95").NumberFormat = "@"
$ws.Range("B95").Value = "4"
$ws.Range("B95").Style = "Normal"
$ws.Range("C95").Value = "happy, exciting"
$ws.Range("D95").Value = "APC"
$ws.Range("E95").Value = "OTH"
$ws.Range("F95").Value = "e9624372-e81d-40ef-b27a-4327fdc73888"
$ws.Range("G95").Value = "BkN_r2lR-_annotated.xlsx"
$ws.Range("H95").Value = "We are also happy that AnonReviewer2 found the list of possible applications, provided in reply to the challenge posted in the review, to be exciting."
